$d = $word.ActiveDocument

# M2DocHTMLParser / HtmlSerializer now always emit explicit bold/italic/strike
# modifiers for the "paragraph style" headings generated from HTML. Mirror
# that by turning on Bold and explicitly turning off Italic and StrikeThrough
# on each of the affected heading runs.

function Set-HeadingModifiers([string]$headingText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($headingText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Font.Bold = $true
        $rng.Font.Italic = $false
        $rng.Font.StrikeThrough = $false
    }
}

Set-HeadingModifiers("Starting with M2Doc")
Set-HeadingModifiers("Template user")
Set-HeadingModifiers("Template developper")
Set-HeadingModifiers("Integrator")
